$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before N (Litter bag 4 area) for the "Blackthorn" species.
# This shifts the old columns N:T to O:U and expands the J1:N1 merge to J1:O1
# (Excel handles this automatically as part of the insert).
$ws.Columns("N").Insert()

# Populate the new "Blackthorn" species sub-header and its data column.
$ws.Range("N2").Value = "Blackthorn"
$ws.Range("N3").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("N6").Value = 0

# Give the new column a custom width, matching how the author resized it.
$ws.Columns("N").ColumnWidth = 9.43

# Add a new data row (row 7) for the period 42306 - 42313.
# Copy the formatting from row 6 first so date styles / blank styles match.
$ws.Range("A6:U6").Copy()
$ws.Range("A7:U7").PasteSpecial(-4122)

$ws.Range("A7").Value = 42306
$ws.Range("B7").Value = 42313

$row7 = @{
    "C" = 12.8
    "D" = 21.47
    "E" = 17.86
    "F" = 2.56
    "G" = 1.62
    "H" = 0
    "I" = 0.88
    "J" = 13.53
    "K" = 9.67
    "L" = 0.42
    "M" = 0.34
    "N" = 0.41
    "O" = 0
    "P" = 10.27
    "Q" = 0.76
    "R" = 1.17
    "S" = 14.91
    "T" = 0
    "U" = 0
}
foreach ($col in $row7.Keys) {
    $ws.Range($col + "7").Value = $row7[$col]
}

# Update the remembered selection, as left behind by the author.
$ws.Range("F11").Select()
